$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Spent - Not from funding" breakdown table in M9:N13
$ws.Range("M9:N13").ClearContents()

# Fill in cost data for existing rows 7-9 (items already listed, previously blank)
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 64.99
$ws.Range("D7").Value = 64.99
$ws.Range("E7").Value = 19.74
$ws.Range("F7").Value = 84.73

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 36.09
$ws.Range("D8").Value = 36.09
$ws.Range("E8").Value = 9.59
$ws.Range("F8").Value = 45.68

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 69.99
$ws.Range("D9").Value = 69.99
$ws.Range("E9").Value = 20.39
$ws.Range("F9").Value = 90.38

$ws.Range("C10").Value = 54.99
$ws.Range("D10").Value = 109.88
$ws.Range("E10").Value = 14.3
$ws.Range("F10").Value = 124.28

# Row 11 already has its values; unchanged.

# New row 12: Aeotec Zwave Smart Switch
$ws.Range("A12").Value = "Aeotec Zwave Smart Switch"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 69.99
$ws.Range("D12").Value = 69.99
$ws.Range("F12").Value = 141.22
$ws.Range("F12").Style = "Accent2"

# New row 13: Everspring Compact Motion Sensor
$ws.Range("A13").Value = "Everspring Compact Motion Sensor"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 44.99
$ws.Range("D13").Value = 44.99
$ws.Range("F13").Style = "Accent2"

# Update selection to reflect the author's final cursor position
$ws.Range("H21").Select()
